$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.19146956285518
$ws.Cells.Item(2, 4).Value = 5.200070649924733
$ws.Cells.Item(2, 5).Value = 18.90384565212386
$ws.Cells.Item(2, 6).Value = 25.81361866815054
$ws.Cells.Item(2, 7).Value = 3.640606607126354
$ws.Cells.Item(2, 11).Value = 9.763602139254225
$ws.Cells.Item(2, 12).Value = 8.41449695919924
$ws.Cells.Item(2, 13).Value = 15.08002433575798
$ws.Cells.Item(2, 14).Value = 20.88808240955116
$ws.Cells.Item(2, 15).Value = 23.09382578303357

$ws.Cells.Item(3, 2).Value = 16.1192693786674
$ws.Cells.Item(3, 4).Value = 5.153896498792251
$ws.Cells.Item(3, 5).Value = 18.95494701023011
$ws.Cells.Item(3, 6).Value = 25.80702669859379
$ws.Cells.Item(3, 7).Value = 3.642414409830277
$ws.Cells.Item(3, 11).Value = 9.454748024583395
$ws.Cells.Item(3, 12).Value = 8.399384925314486
$ws.Cells.Item(3, 13).Value = 15.0667106688066
$ws.Cells.Item(3, 14).Value = 20.95097965815719
$ws.Cells.Item(3, 15).Value = 23.1370932920056

$ws.Cells.Item(4, 2).Value = 16.07791005306803
$ws.Cells.Item(4, 4).Value = 5.124844573883039
$ws.Cells.Item(4, 5).Value = 18.98811724536805
$ws.Cells.Item(4, 6).Value = 25.80946186041148
$ws.Cells.Item(4, 7).Value = 3.643584334350189
$ws.Cells.Item(4, 11).Value = 9.257953156623193
$ws.Cells.Item(4, 12).Value = 8.391308362560068
$ws.Cells.Item(4, 13).Value = 15.06072546774611
$ws.Cells.Item(4, 14).Value = 20.9914256224606
$ws.Cells.Item(4, 15).Value = 23.16847341401841

$ws.Cells.Item(5, 2).Value = 16.06181697555458
$ws.Cells.Item(5, 4).Value = 5.112833913498615
$ws.Cells.Item(5, 5).Value = 19.00208657253476
$ws.Cells.Item(5, 6).Value = 25.81208574776883
$ws.Cells.Item(5, 7).Value = 3.644076203566508
$ws.Cells.Item(5, 11).Value = 9.176028475885238
$ws.Cells.Item(5, 12).Value = 8.38832191083662
$ws.Cells.Item(5, 13).Value = 15.05883992322635
$ws.Cells.Item(5, 14).Value = 21.00836838530107
$ws.Cells.Item(5, 15).Value = 23.18246977829369

$ws.Cells.Item(6, 2).Value = 16.05919109287119
$ws.Cells.Item(6, 4).Value = 5.110829309813401
$ws.Cells.Item(6, 5).Value = 19.00443351274773
$ws.Cells.Item(6, 6).Value = 25.81262000345958
$ws.Cells.Item(6, 7).Value = 3.644158792381177
$ws.Cells.Item(6, 11).Value = 9.162322707394436
$ws.Cells.Item(6, 12).Value = 8.387844491612114
$ws.Cells.Item(6, 13).Value = 15.05856033852333
$ws.Cells.Item(6, 14).Value = 21.01120958299605
$ws.Cells.Item(6, 15).Value = 23.18486681107413

$ws.Cells.Item(7, 2).Value = 16.07768991656132
$ws.Cells.Item(7, 4).Value = 5.124683283459479
$ws.Cells.Item(7, 5).Value = 18.9883038082053
$ws.Cells.Item(7, 6).Value = 25.80949063997614
$ws.Cells.Item(7, 7).Value = 3.643590906606336
$ws.Cells.Item(7, 11).Value = 9.256855194703602
$ws.Cells.Item(7, 12).Value = 8.391266848908968
$ws.Cells.Item(7, 13).Value = 15.06069779397976
$ws.Cells.Item(7, 14).Value = 20.99165225133617
$ws.Cells.Item(7, 15).Value = 23.16865728236547

$ws.Cells.Item(8, 2).Value = 16.16596673489348
$ws.Cells.Item(8, 4).Value = 5.184297906480111
$ws.Cells.Item(8, 5).Value = 18.9210937748824
$ws.Cells.Item(8, 6).Value = 25.81000145351296
$ws.Cells.Item(8, 7).Value = 3.641217527817932
$ws.Cells.Item(8, 11).Value = 9.65863488544744
$ws.Cells.Item(8, 12).Value = 8.409038344388522
$ws.Cells.Item(8, 13).Value = 15.07498119149762
$ws.Cells.Item(8, 14).Value = 20.90939103918387
$ws.Cells.Item(8, 15).Value = 23.10774437212165

$ws.Cells.Item(9, 2).Value = 16.36200606840484
$ws.Cells.Item(9, 4).Value = 5.295450714934856
$ws.Cells.Item(9, 5).Value = 18.80347691183033
$ws.Cells.Item(9, 6).Value = 25.86233217878239
$ws.Cells.Item(9, 7).Value = 3.637036707511726
$ws.Cells.Item(9, 11).Value = 10.38698588073074
$ws.Cells.Item(9, 12).Value = 8.45331094983214
$ws.Cells.Item(9, 13).Value = 15.12022521344731
$ws.Cells.Item(9, 14).Value = 20.76250970290297
$ws.Cells.Item(9, 15).Value = 23.02655678634418

$ws.Cells.Item(10, 2).Value = 16.51908405129227
$ws.Cells.Item(10, 4).Value = 5.373361001665589
$ws.Cells.Item(10, 5).Value = 18.72563794956697
$ws.Cells.Item(10, 6).Value = 25.93185921044816
$ws.Cells.Item(10, 7).Value = 3.634250671245182
$ws.Cells.Item(10, 11).Value = 10.88259205856728
$ws.Cells.Item(10, 12).Value = 8.49141494325659
$ws.Cells.Item(10, 13).Value = 15.16377344061464
$ws.Cells.Item(10, 14).Value = 20.66330602861025
$ws.Cells.Item(10, 15).Value = 22.9903124978443

$ws.Cells.Item(11, 2).Value = 16.5931617685295
$ws.Cells.Item(11, 4).Value = 5.407938622775498
$ws.Cells.Item(11, 5).Value = 18.69207365664418
$ws.Cells.Item(11, 6).Value = 25.97016870301998
$ws.Cells.Item(11, 7).Value = 3.633044619025731
$ws.Cells.Item(11, 11).Value = 11.0989171440531
$ws.Cells.Item(11, 12).Value = 8.50992083902355
$ws.Cells.Item(11, 13).Value = 15.18577493115219
$ws.Cells.Item(11, 14).Value = 20.62004812797257
$ws.Cells.Item(11, 15).Value = 22.97891669185882

$ws.Cells.Item(12, 2).Value = 16.62157102375949
$ws.Cells.Item(12, 4).Value = 5.420903811181253
$ws.Cells.Item(12, 5).Value = 18.67962791398938
$ws.Cells.Item(12, 6).Value = 25.98562852471273
$ws.Cells.Item(12, 7).Value = 3.632596689834988
$ws.Cells.Item(12, 11).Value = 11.17948157950919
$ws.Cells.Item(12, 12).Value = 8.517093434586776
$ws.Cells.Item(12, 13).Value = 15.19441688257587
$ws.Cells.Item(12, 14).Value = 20.60393504627305
$ws.Cells.Item(12, 15).Value = 22.97533388457427

$ws.Cells.Item(13, 2).Value = 16.61543699201549
$ws.Cells.Item(13, 4).Value = 5.418117313602917
$ws.Cells.Item(13, 5).Value = 18.68229658697339
$ws.Cells.Item(13, 6).Value = 25.98225673558009
$ws.Cells.Item(13, 7).Value = 3.632692769731842
$ws.Cells.Item(13, 11).Value = 11.16219136224831
$ws.Cells.Item(13, 12).Value = 8.515541419155296
$ws.Cells.Item(13, 13).Value = 15.19254195218016
$ws.Cells.Item(13, 14).Value = 20.60739339772694
$ws.Cells.Item(13, 15).Value = 22.97607292490708

$ws.Cells.Item(14, 2).Value = 16.59549196091463
$ws.Cells.Item(14, 4).Value = 5.409007878556652
$ws.Cells.Item(14, 5).Value = 18.69104444602526
$ws.Cells.Item(14, 6).Value = 25.97142154478276
$ws.Cells.Item(14, 7).Value = 3.63300759196802
$ws.Cells.Item(14, 11).Value = 11.10557256276057
$ws.Cells.Item(14, 12).Value = 8.510507651734706
$ws.Cells.Item(14, 13).Value = 15.1864797149589
$ws.Cells.Item(14, 14).Value = 20.61871713588191
$ws.Cells.Item(14, 15).Value = 22.97860725036532

$ws.Cells.Item(15, 2).Value = 16.58332104114573
$ws.Cells.Item(15, 4).Value = 5.403411211299295
$ws.Cells.Item(15, 5).Value = 18.6964371562834
$ws.Cells.Item(15, 6).Value = 25.9649085119721
$ws.Cells.Item(15, 7).Value = 3.633201571395438
$ws.Cells.Item(15, 11).Value = 11.07071457214137
$ws.Cells.Item(15, 12).Value = 8.507445673350917
$ws.Cells.Item(15, 13).Value = 15.18280670946652
$ws.Cells.Item(15, 14).Value = 20.62568808531844
$ws.Cells.Item(15, 15).Value = 22.98025499966763

$ws.Cells.Item(16, 2).Value = 16.51429401506332
$ws.Cells.Item(16, 4).Value = 5.371083510077336
$ws.Cells.Item(16, 5).Value = 18.72786849727683
$ws.Cells.Item(16, 6).Value = 25.92948938966939
$ws.Cells.Item(16, 7).Value = 3.634330720059896
$ws.Cells.Item(16, 11).Value = 10.86826743817062
$ws.Cells.Item(16, 12).Value = 8.490228810904595
$ws.Cells.Item(16, 13).Value = 15.16237929439477
$ws.Cells.Item(16, 14).Value = 20.66617056687341
$ws.Cells.Item(16, 15).Value = 22.99115973854891

$ws.Cells.Item(17, 2).Value = 16.47260613106224
$ws.Cells.Item(17, 4).Value = 5.351027056802197
$ws.Cells.Item(17, 5).Value = 18.74762245966821
$ws.Cells.Item(17, 6).Value = 25.90946688266458
$ws.Cells.Item(17, 7).Value = 3.635039093819835
$ws.Cells.Item(17, 11).Value = 10.7417042003425
$ws.Cells.Item(17, 12).Value = 8.479964442997739
$ws.Cells.Item(17, 13).Value = 15.15040574561467
$ws.Cells.Item(17, 14).Value = 20.69148343652658
$ws.Cells.Item(17, 15).Value = 22.99915394340257

$ws.Cells.Item(18, 2).Value = 16.44887637518998
$ws.Cells.Item(18, 4).Value = 5.33941002941462
$ws.Cells.Item(18, 5).Value = 18.75915813878775
$ws.Cells.Item(18, 6).Value = 25.89857988902015
$ws.Cells.Item(18, 7).Value = 3.635452306826607
$ws.Cells.Item(18, 11).Value = 10.66805122887423
$ws.Cells.Item(18, 12).Value = 8.474171218612979
$ws.Cells.Item(18, 13).Value = 15.14372541544216
$ws.Cells.Item(18, 14).Value = 20.70621884713101
$ws.Cells.Item(18, 15).Value = 23.00423127455234

$ws.Cells.Item(19, 2).Value = 16.44088506632865
$ws.Cells.Item(19, 4).Value = 5.335462911056541
$ws.Cells.Item(19, 5).Value = 18.76309378880483
$ws.Cells.Item(19, 6).Value = 25.89500206621425
$ws.Cells.Item(19, 7).Value = 3.635593206841336
$ws.Cells.Item(19, 11).Value = 10.64296769897024
$ws.Cells.Item(19, 12).Value = 8.472228833848634
$ws.Cells.Item(19, 13).Value = 15.14149918015596
$ws.Cells.Item(19, 14).Value = 20.71123828930336
$ws.Cells.Item(19, 15).Value = 23.006032672145

$ws.Cells.Item(20, 2).Value = 16.47701834803281
$ws.Cells.Item(20, 4).Value = 5.353170524376416
$ws.Cells.Item(20, 5).Value = 18.74550164423152
$ws.Cells.Item(20, 6).Value = 25.91153322346473
$ws.Cells.Item(20, 7).Value = 3.634963088825947
$ws.Cells.Item(20, 11).Value = 10.75526609341066
$ws.Cells.Item(20, 12).Value = 8.481045685626173
$ws.Cells.Item(20, 13).Value = 15.15165900703346
$ws.Cells.Item(20, 14).Value = 20.68877061961499
$ws.Cells.Item(20, 15).Value = 22.99825333990352

$ws.Cells.Item(21, 2).Value = 16.60134075169242
$ws.Cells.Item(21, 4).Value = 5.411687064822138
$ws.Cells.Item(21, 5).Value = 18.68846782184118
$ws.Cells.Item(21, 6).Value = 25.97457831172311
$ws.Cells.Item(21, 7).Value = 3.632914883187818
$ws.Cells.Item(21, 11).Value = 11.12223989974962
$ws.Cells.Item(21, 12).Value = 8.511981748959611
$ws.Cells.Item(21, 13).Value = 15.18825195323688
$ws.Cells.Item(21, 14).Value = 20.61538382249851
$ws.Cells.Item(21, 15).Value = 22.97784297580132

$ws.Cells.Item(22, 2).Value = 16.68466778111519
$ws.Cells.Item(22, 4).Value = 5.449179457824598
$ws.Cells.Item(22, 5).Value = 18.65273316926706
$ws.Cells.Item(22, 6).Value = 26.02133137026592
$ws.Cells.Item(22, 7).Value = 3.631627399721854
$ws.Cells.Item(22, 11).Value = 11.35417574567986
$ws.Cells.Item(22, 12).Value = 8.533159059042095
$ws.Cells.Item(22, 13).Value = 15.21397485271923
$ws.Cells.Item(22, 14).Value = 20.56898138718966
$ws.Cells.Item(22, 15).Value = 22.96877331618133

$ws.Cells.Item(23, 2).Value = 16.64001126341718
$ws.Cells.Item(23, 4).Value = 5.429239233705319
$ws.Cells.Item(23, 5).Value = 18.67166481666725
$ws.Cells.Item(23, 6).Value = 25.99587351937807
$ws.Cells.Item(23, 7).Value = 3.632309888765178
$ws.Cells.Item(23, 11).Value = 11.23112244183626
$ws.Cells.Item(23, 12).Value = 8.521769882231915
$ws.Cells.Item(23, 13).Value = 15.20008228563172
$ws.Cells.Item(23, 14).Value = 20.59360488888441
$ws.Cells.Item(23, 15).Value = 22.97322326848711

$ws.Cells.Item(24, 2).Value = 16.47502284453081
$ws.Cells.Item(24, 4).Value = 5.352201731175917
$ws.Cells.Item(24, 5).Value = 18.74645990768968
$ws.Cells.Item(24, 6).Value = 25.9105970858199
$ws.Cells.Item(24, 7).Value = 3.634997432116851
$ws.Cells.Item(24, 11).Value = 10.74913753087385
$ws.Cells.Item(24, 12).Value = 8.480556519464203
$ws.Cells.Item(24, 13).Value = 15.15109177355949
$ws.Cells.Item(24, 14).Value = 20.68999651492484
$ws.Cells.Item(24, 15).Value = 22.99865900336913

$ws.Cells.Item(25, 2).Value = 16.30661032708882
$ws.Cells.Item(25, 4).Value = 5.266022531005729
$ws.Cells.Item(25, 5).Value = 18.83378455802354
$ws.Cells.Item(25, 6).Value = 25.84269752608531
$ws.Cells.Item(25, 7).Value = 3.638117360729984
$ws.Cells.Item(25, 11).Value = 10.38698588073074
$ws.Cells.Item(25, 12).Value = 8.440341776209713
$ws.Cells.Item(25, 13).Value = 15.12022521344731
$ws.Cells.Item(25, 14).Value = 20.76250970290297
$ws.Cells.Item(25, 15).Value = 23.02655678634418

Write-Host "Applied loading_percent updates for case with 380 kV"
